# Update the division-problem table with the newly generated values.
# Each source row (1,5,9,13,17 -- the rows that actually contain text)
# has 5 cells; find/replace is scoped to each individual cell Range so
# that duplicate "old" values in different cells map to the correct,
# distinct "new" values.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$failures = 0

$r = $t.Cell(1,1).Range
$ok = $r.Find.Execute("90÷9=10, 0", $true, $false, $false, $false, $false, $true, 0, $false, "78÷9=8, 6", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(1,1): 90÷9=10, 0 -> 78÷9=8, 6" }

$r = $t.Cell(1,2).Range
$ok = $r.Find.Execute("88÷8=11, 0", $true, $false, $false, $false, $false, $true, 0, $false, "50÷4=12, 2", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(1,2): 88÷8=11, 0 -> 50÷4=12, 2" }

$r = $t.Cell(1,3).Range
$ok = $r.Find.Execute("77÷5=15, 2", $true, $false, $false, $false, $false, $true, 0, $false, "48÷3=16, 0", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(1,3): 77÷5=15, 2 -> 48÷3=16, 0" }

$r = $t.Cell(1,4).Range
$ok = $r.Find.Execute("82÷2=41, 0", $true, $false, $false, $false, $false, $true, 0, $false, "48÷8=6, 0", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(1,4): 82÷2=41, 0 -> 48÷8=6, 0" }

$r = $t.Cell(1,5).Range
$ok = $r.Find.Execute("60÷8=7, 4", $true, $false, $false, $false, $false, $true, 0, $false, "41÷6=6, 5", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(1,5): 60÷8=7, 4 -> 41÷6=6, 5" }

$r = $t.Cell(5,1).Range
$ok = $r.Find.Execute("48÷9=5, 3", $true, $false, $false, $false, $false, $true, 0, $false, "60÷7=8, 4", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(5,1): 48÷9=5, 3 -> 60÷7=8, 4" }

$r = $t.Cell(5,2).Range
$ok = $r.Find.Execute("40÷8=5, 0", $true, $false, $false, $false, $false, $true, 0, $false, "56÷3=18, 2", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(5,2): 40÷8=5, 0 -> 56÷3=18, 2" }

$r = $t.Cell(5,3).Range
$ok = $r.Find.Execute("29÷5=5, 4", $true, $false, $false, $false, $false, $true, 0, $false, "93÷9=10, 3", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(5,3): 29÷5=5, 4 -> 93÷9=10, 3" }

$r = $t.Cell(5,4).Range
$ok = $r.Find.Execute("91÷6=15, 1", $true, $false, $false, $false, $false, $true, 0, $false, "25÷2=12, 1", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(5,4): 91÷6=15, 1 -> 25÷2=12, 1" }

$r = $t.Cell(5,5).Range
$ok = $r.Find.Execute("99÷9=11, 0", $true, $false, $false, $false, $false, $true, 0, $false, "77÷5=15, 2", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(5,5): 99÷9=11, 0 -> 77÷5=15, 2" }

$r = $t.Cell(9,1).Range
$ok = $r.Find.Execute("41÷9=4, 5", $true, $false, $false, $false, $false, $true, 0, $false, "97÷6=16, 1", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(9,1): 41÷9=4, 5 -> 97÷6=16, 1" }

$r = $t.Cell(9,2).Range
$ok = $r.Find.Execute("53÷6=8, 5", $true, $false, $false, $false, $false, $true, 0, $false, "77÷2=38, 1", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(9,2): 53÷6=8, 5 -> 77÷2=38, 1" }

$r = $t.Cell(9,3).Range
$ok = $r.Find.Execute("29÷9=3, 2", $true, $false, $false, $false, $false, $true, 0, $false, "41÷7=5, 6", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(9,3): 29÷9=3, 2 -> 41÷7=5, 6" }

$r = $t.Cell(9,4).Range
$ok = $r.Find.Execute("85÷7=12, 1", $true, $false, $false, $false, $false, $true, 0, $false, "14÷2=7, 0", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(9,4): 85÷7=12, 1 -> 14÷2=7, 0" }

$r = $t.Cell(9,5).Range
$ok = $r.Find.Execute("90÷9=10, 0", $true, $false, $false, $false, $false, $true, 0, $false, "74÷6=12, 2", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(9,5): 90÷9=10, 0 -> 74÷6=12, 2" }

$r = $t.Cell(13,1).Range
$ok = $r.Find.Execute("28÷7=4, 0", $true, $false, $false, $false, $false, $true, 0, $false, "10÷8=1, 2", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(13,1): 28÷7=4, 0 -> 10÷8=1, 2" }

$r = $t.Cell(13,2).Range
$ok = $r.Find.Execute("67÷8=8, 3", $true, $false, $false, $false, $false, $true, 0, $false, "81÷5=16, 1", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(13,2): 67÷8=8, 3 -> 81÷5=16, 1" }

$r = $t.Cell(13,3).Range
$ok = $r.Find.Execute("97÷7=13, 6", $true, $false, $false, $false, $false, $true, 0, $false, "76÷4=19, 0", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(13,3): 97÷7=13, 6 -> 76÷4=19, 0" }

$r = $t.Cell(13,4).Range
$ok = $r.Find.Execute("23÷5=4, 3", $true, $false, $false, $false, $false, $true, 0, $false, "30÷3=10, 0", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(13,4): 23÷5=4, 3 -> 30÷3=10, 0" }

$r = $t.Cell(13,5).Range
$ok = $r.Find.Execute("17÷4=4, 1", $true, $false, $false, $false, $false, $true, 0, $false, "91÷3=30, 1", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(13,5): 17÷4=4, 1 -> 91÷3=30, 1" }

$r = $t.Cell(17,1).Range
$ok = $r.Find.Execute("97÷2=48, 1", $true, $false, $false, $false, $false, $true, 0, $false, "20÷9=2, 2", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(17,1): 97÷2=48, 1 -> 20÷9=2, 2" }

$r = $t.Cell(17,2).Range
$ok = $r.Find.Execute("64÷8=8, 0", $true, $false, $false, $false, $false, $true, 0, $false, "96÷4=24, 0", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(17,2): 64÷8=8, 0 -> 96÷4=24, 0" }

$r = $t.Cell(17,3).Range
$ok = $r.Find.Execute("69÷5=13, 4", $true, $false, $false, $false, $false, $true, 0, $false, "46÷7=6, 4", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(17,3): 69÷5=13, 4 -> 46÷7=6, 4" }

$r = $t.Cell(17,4).Range
$ok = $r.Find.Execute("88÷6=14, 4", $true, $false, $false, $false, $false, $true, 0, $false, "34÷8=4, 2", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(17,4): 88÷6=14, 4 -> 34÷8=4, 2" }

$r = $t.Cell(17,5).Range
$ok = $r.Find.Execute("82÷2=41, 0", $true, $false, $false, $false, $false, $true, 0, $false, "89÷6=14, 5", 1)
if (-not $ok) { $failures++; Write-Output "FAILED Cell(17,5): 82÷2=41, 0 -> 89÷6=14, 5" }

Write-Output ("Done. Failures=" + $failures)